# EPI-634 Make vital field mandatory by encounter type (#4762)
#
# Adds two new "Succeed" rows to the "Question Validation Succeed" sheet
# covering validationCriteria.mandatory keyed by encounterType (single
# value and array-of-values forms), widens column A to fit the new,
# longer row labels, and nudges the saved selections/page setup to match
# the author's final state.

$wb = $excel.ActiveWorkbook
$wsSucceed = $wb.Worksheets.Item("Question Validation Succeed")
$wsFail = $wb.Worksheets.Item("Question Validation Fail")

# --- New data rows -------------------------------------------------------
# Row 12's label is entered before row 11's so that the shared-string
# table ends up in the same order as the source workbook.
$wsSucceed.Range("A12").Value = "Succeed-validationCriteria-mandatory-2"
$wsSucceed.Range("A11").Value = "Succeed-validationCriteria-mandatory-1"
$wsSucceed.Range("K11").Value = '{ "mandatory": {"encounterType":"admission"} }'
$wsSucceed.Range("K12").Value = '{ "mandatory": {"encounterType":["admission","surveyResponse"]} }'

$wsSucceed.Range("B11").Value = "SurveyAnswer"
$wsSucceed.Range("C11").Value = "SurveyAnswer: Full config"
$wsSucceed.Range("P11").Value = '{ "source": "xyz" }'

$wsSucceed.Range("B12").Value = "SurveyAnswer"
$wsSucceed.Range("C12").Value = "SurveyAnswer: Full config"
$wsSucceed.Range("P12").Value = '{ "source": "xyz" }'

# --- Column widths ---------------------------------------------------
# Column A needs to widen (and drop its autofit/bestFit flag) to fit the
# new, longer row labels.
$wsSucceed.Columns.Item(1).ColumnWidth = 56

# --- Page setup --------------------------------------------------------
$wsSucceed.PageSetup.PaperSize = 9
$wsSucceed.PageSetup.Orientation = 1

# --- Selections ----------------------------------------------------------
# Update the "Question Validation Fail" sheet's remembered selection
# without leaving it as the active tab.
$wsFail.Range("K12").Select()
$wsSucceed.Activate()
$wsSucceed.Range("P17").Select()
